$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 910349.8
$ws.Range("I2").Value = 524.75
$ws.Range("J2").Value = 1430249.9
$ws.Range("K2").Value = 524.75
$ws.Range("L2").Value = 1430249.9
$ws.Range("M2").Value = -411.75
$ws.Range("N2").Value = -1430475.9
$ws.Range("H4").Value = 1087.381
$ws.Range("I4").Value = 802.0714
$ws.Range("K4").Value = 802.0714
$ws.Range("M4").Value = -688.0714
$ws.Range("H18").Value = 2598.5
$ws.Range("I18").Value = 2598.5
$ws.Range("J18").Value = 0.0
$ws.Range("K18").Value = 2598.5
$ws.Range("L18").Value = 0.0
$ws.Range("M18").Value = -2314.5
$ws.Range("N18").ClearContents()
$ws.Range("H33").Value = 311.33334
$ws.Range("I33").Value = 282.2
$ws.Range("K33").Value = 282.2
$ws.Range("M33").Value = -53.19999999999999
$ws.Range("H38").Value = 2099.5557
$ws.Range("I38").Value = 128.0
$ws.Range("J38").Value = 9000.0
$ws.Range("K38").Value = 384.0
$ws.Range("L38").Value = 27000.0
$ws.Range("M38").Value = -12.0
$ws.Range("N38").Value = -27744.0
$ws.Range("H39").Value = 225.0
$ws.Range("I39").Value = 103.666664
$ws.Range("J39").Value = 407.0
$ws.Range("K39").Value = 310.999992
$ws.Range("L39").Value = 1221.0
$ws.Range("M39").Value = -14.99999200000002
$ws.Range("N39").Value = -1813.0
$ws.Range("H40").Value = 50002988.0
$ws.Range("J40").Value = 125002376.0
$ws.Range("L40").Value = 125002376.0
$ws.Range("N40").Value = -125002726.0
$ws.Range("H42").Value = 189.22223
$ws.Range("I42").Value = 140.9
$ws.Range("J42").Value = 249.625
$ws.Range("K42").Value = 422.7
$ws.Range("L42").Value = 748.875
$ws.Range("M42").Value = -192.7
$ws.Range("N42").Value = -1208.875
$ws.Range("H48").Value = 0.0
$ws.Range("J48").Value = 0.0
$ws.Range("L48").Value = 0.0
$ws.Range("N48").ClearContents()
$ws.Range("H51").Value = 7832.857
$ws.Range("J51").Value = 5736.8423
$ws.Range("L51").Value = 5736.8423
$ws.Range("N51").Value = -6704.8423
$ws.Range("H52").Value = 1050.0
$ws.Range("I52").Value = 1050.0
$ws.Range("J52").Value = 0.0
$ws.Range("K52").Value = 3150.0
$ws.Range("L52").Value = 0.0
$ws.Range("M52").Value = -2990.0
$ws.Range("N52").ClearContents()
$ws.Range("H56").Value = 0.0
$ws.Range("J56").Value = 0.0
$ws.Range("L56").Value = 0.0
$ws.Range("N56").ClearContents()
$ws.Range("H58").Value = 6609.0
$ws.Range("I58").Value = 327.5
$ws.Range("J58").Value = 9749.75
$ws.Range("K58").Value = 982.5
$ws.Range("L58").Value = 29249.25
$ws.Range("M58").Value = -832.5
$ws.Range("N58").Value = -29549.25
$ws.Range("H64").Value = 3647.3333
$ws.Range("I64").Value = 3501.0
$ws.Range("K64").Value = 3501.0
$ws.Range("M64").Value = -3253.0
$ws.Range("H67").Value = 3647.3333
$ws.Range("I67").Value = 3501.0
$ws.Range("K67").Value = 3501.0
$ws.Range("M67").Value = -2643.0
$ws.Range("H106").Value = 13968.5
$ws.Range("I106").Value = 12424.667
$ws.Range("K106").Value = 12424.667
$ws.Range("M106").Value = -11793.667
$ws.Range("H135").Value = 2641.4666
$ws.Range("J135").Value = 4544.3335
$ws.Range("L135").Value = 40899.0015
$ws.Range("N135").Value = -45969.0015
$ws.Range("H137").Value = 3205.3333
$ws.Range("I137").Value = 2668.5
$ws.Range("K137").Value = 8005.5
$ws.Range("M137").Value = -5455.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 125.833336
$ws.Range("H23").Value = 7500.0
$ws.Range("J23").Value = 7500.0
$ws.Range("L23").Value = 7500.0
$ws.Range("N23").Value = -8018.0
$ws.Range("H32").Value = 4556.4424
$ws.Range("I32").Value = 3957.8572
$ws.Range("K32").Value = 3957.8572
$ws.Range("M32").Value = -3670.8572
$ws.Range("H61").Value = 11581445.0
$ws.Range("I61").Value = 13127717.0
$ws.Range("K61").Value = 13127717.0
$ws.Range("M61").Value = -13127505.0
$ws.Range("H88").Value = 2236.5557
$ws.Range("I88").Value = 1465.6666
$ws.Range("K88").Value = 1465.6666
$ws.Range("M88").Value = -1059.6666
$ws.Range("H91").Value = 2236.5557
$ws.Range("I91").Value = 1465.6666
$ws.Range("K91").Value = 1465.6666
$ws.Range("M91").Value = -61.66660000000002
$ws.Range("H110").Value = 6251.9
$ws.Range("I110").Value = 7589.4287
$ws.Range("K110").Value = 7589.4287
$ws.Range("M110").Value = -5544.4287
$ws.Range("H136").Value = 11581445.0
$ws.Range("I136").Value = 13127717.0
$ws.Range("K136").Value = 39383151.0
$ws.Range("M136").Value = -39380601.0

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 125.833336
$ws.Range("H18").Value = 9998.5
$ws.Range("I18").Value = 9998.0
$ws.Range("J18").Value = 9999.0
$ws.Range("K18").Value = 9998.0
$ws.Range("L18").Value = 9999.0
$ws.Range("M18").Value = -9469.0
$ws.Range("N18").Value = -11057.0
$ws.Range("H19").Value = 9999.0
$ws.Range("J19").Value = 9999.0
$ws.Range("L19").Value = 9999.0
$ws.Range("N19").Value = -10345.0
$ws.Range("H22").Value = 639.6667
$ws.Range("I22").Value = 720.0
$ws.Range("J22").Value = 238.0
$ws.Range("K22").Value = 720.0
$ws.Range("L22").Value = 238.0
$ws.Range("M22").Value = -547.0
$ws.Range("N22").Value = -584.0
$ws.Range("H99").Value = 1816.0
$ws.Range("I99").Value = 1855.5625
$ws.Range("K99").Value = 1855.5625
$ws.Range("M99").Value = -357.5625
$ws.Range("I134").Value = 2456.25
$ws.Range("K134").Value = 7368.75
$ws.Range("M134").Value = -4833.75
$ws.Range("H135").Value = 100000.0
$ws.Range("J135").Value = 100000.0
$ws.Range("L135").Value = 100000.0
$ws.Range("N135").Value = -110140.0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 161.0
$ws.Range("I7").Value = 90.0
$ws.Range("J7").Value = 800.0
$ws.Range("K7").Value = 90.0
$ws.Range("L7").Value = 800.0
$ws.Range("M7").Value = 23.0
$ws.Range("N7").Value = -1026.0
$ws.Range("H16").Value = 3450524.8
$ws.Range("I16").Value = 3848381.5
$ws.Range("K16").Value = 3848381.5
$ws.Range("M16").Value = -3848094.5
$ws.Range("H113").Value = 3450524.8
$ws.Range("I113").Value = 3848381.5
$ws.Range("K113").Value = 3848381.5
$ws.Range("M113").Value = -3846211.5
$ws.Range("H122").Value = 4322.4375
$ws.Range("I122").Value = 4302.222
$ws.Range("K122").Value = 12906.666
$ws.Range("M122").Value = -10456.666

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1350.1666
$ws.Range("I5").Value = 1260.2
$ws.Range("K5").Value = 3780.6
$ws.Range("M5").Value = -3668.6
$ws.Range("H23").Value = 501.5
$ws.Range("I23").Value = 66.0
$ws.Range("J23").Value = 719.25
$ws.Range("K23").Value = 198.0
$ws.Range("L23").Value = 2157.75
$ws.Range("M23").Value = 37.0
$ws.Range("N23").Value = -2627.75
$ws.Range("H24").Value = 16777.5
$ws.Range("I24").Value = 0.0
$ws.Range("K24").Value = 0.0
$ws.Range("M24").ClearContents()
$ws.Range("H25").Value = 14065.8
$ws.Range("I25").Value = 8999.667
$ws.Range("K25").Value = 26999.001
$ws.Range("M25").Value = -26830.001
$ws.Range("H30").Value = 14065.8
$ws.Range("I30").Value = 8999.667
$ws.Range("K30").Value = 26999.001
$ws.Range("M30").Value = -26897.001
$ws.Range("H135").Value = 1350.1666
$ws.Range("I135").Value = 1260.2
$ws.Range("K135").Value = 11341.8
$ws.Range("M135").Value = -8806.800000000001

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2629.7222
$ws.Range("I126").Value = 1971.2307
$ws.Range("J126").Value = 4341.8
$ws.Range("K126").Value = 5913.6921
$ws.Range("L126").Value = 13025.4
$ws.Range("M126").Value = -3443.6921
$ws.Range("N126").Value = -17965.4
$ws.Range("H132").Value = 5557760.0
$ws.Range("I132").Value = 1910.1666
$ws.Range("J132").Value = 16669460.0
$ws.Range("K132").Value = 5730.4998
$ws.Range("L132").Value = 50008380.0
$ws.Range("M132").Value = -3200.4998
$ws.Range("N132").Value = -50013440.0
$ws.Range("H135").Value = 107910.25
$ws.Range("J135").Value = 107910.25
$ws.Range("L135").Value = 107910.25
$ws.Range("N135").Value = -118050.25

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4426.9473
$ws.Range("I16").Value = 1838.0769
$ws.Range("J16").Value = 10036.167
$ws.Range("K16").Value = 1838.0769
$ws.Range("L16").Value = 10036.167
$ws.Range("M16").Value = -1668.0769
$ws.Range("N16").Value = -10376.167
$ws.Range("H61").Value = 1892.4445
$ws.Range("I61").Value = 1005.0
$ws.Range("J61").Value = 4998.5
$ws.Range("K61").Value = 1005.0
$ws.Range("L61").Value = 4998.5
$ws.Range("M61").Value = -803.0
$ws.Range("N61").Value = -5402.5
$ws.Range("H113").Value = 1892.4445
$ws.Range("I113").Value = 1005.0
$ws.Range("J113").Value = 4998.5
$ws.Range("K113").Value = 1005.0
$ws.Range("L113").Value = 4998.5
$ws.Range("M113").Value = 1165.0
$ws.Range("N113").Value = -9338.5
$ws.Range("H122").Value = 3357.1904
$ws.Range("I122").Value = 3026.4211
$ws.Range("K122").Value = 9079.2633
$ws.Range("M122").Value = -6629.263300000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 302768.7
$ws.Range("J136").Value = 1114244.9
$ws.Range("L136").Value = 3342734.7
$ws.Range("N136").Value = -3347834.7
